# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# Update column G ("K") values for rows 2-17 on the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 2
    3  = 1
    4  = 1
    5  = 0
    6  = 2
    7  = 2
    8  = 1
    9  = 1
    10 = 1
    11 = 2
    12 = 2
    13 = 2
    14 = 1
    15 = 2
    16 = 2
    17 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
